# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 163, pushing the existing
# rows 163..252 down to 164..253 (dimension grows from A1:T252 to A1:T253).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 163; formatting (date style on
# column D) carries over automatically from the neighbouring rows.
$ws.Rows.Item(163).Insert()

$ws.Cells.Item(163, 1).Value  = 4
$ws.Cells.Item(163, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(163, 3).Value  = "Los Lagos"
$ws.Cells.Item(163, 4).Value  = 44830
$ws.Cells.Item(163, 5).Value  = 10
$ws.Cells.Item(163, 6).Value  = "Fruta"
$ws.Cells.Item(163, 7).Value  = 100108
$ws.Cells.Item(163, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(163, 9).Value  = 100108002
$ws.Cells.Item(163, 10).Value = "Mango"
$ws.Cells.Item(163, 11).Value = "Sin especificar"
$ws.Cells.Item(163, 12).Value = "Primera"
$ws.Cells.Item(163, 13).Value = 160
$ws.Cells.Item(163, 14).Value = 9000
$ws.Cells.Item(163, 15).Value = 10000
$ws.Cells.Item(163, 16).Value = 9500
$ws.Cells.Item(163, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(163, 18).Value = "Brasil"
$ws.Cells.Item(163, 19).Value = 2375
$ws.Cells.Item(163, 20).Value = 4
